$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 header fields (C and D only changed) ---
# Row 8
$ws.Range("C8").Value = "MLK_PMT_10107_-_H-001"
$ws.Range("D8").Value = "Cooling of Water on Irrigation of An Absorber"
$ws.Range("E8").Value = "Shell"
$ws.Range("G8").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H8").Value = "Carbon Steel"
$ws.Range("I8").Value = "SA-516"
$ws.Range("J8").Value = "Gr. 70"
$ws.Range("L8").Value = "200 DEG. C"
$ws.Range("M8").Value = "14.00 kg/cm2 (g)"
$ws.Range("N8").Value = "150 DEG. C"
$ws.Range("O8").Value = "10.00 kg/cm2 (g)"

# Row 9
$ws.Range("E9").Value = "Shell Cover (FLANGE)"
$ws.Range("G9").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H9").Value = "Carbon Steel"
$ws.Range("I9").Value = "SA-516"
$ws.Range("J9").Value = "Gr. 70"
$ws.Range("L9").Value = "200 DEG. C"
$ws.Range("M9").Value = "14.00 kg/cm2 (g)"
$ws.Range("N9").Value = "150 DEG. C"
$ws.Range("O9").Value = "10.00 kg/cm2 (g)"

# Row 10
$ws.Range("E10").Value = "Tubesheet"
$ws.Range("G10").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H10").Value = "Carbon Steel"
$ws.Range("I10").Value = "SA-516"
$ws.Range("J10").Value = "Gr. 70"
$ws.Range("L10").Value = "200 DEG. C"
$ws.Range("M10").Value = "14.00 kg/cm2 (g)"
$ws.Range("N10").Value = "150 DEG. C"
$ws.Range("O10").Value = "10.00 kg/cm2 (g)"

# Row 11
$ws.Range("E11").Value = "Head"
$ws.Range("G11").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H11").Value = "Carbon Steel"
$ws.Range("I11").Value = "SA-516"
$ws.Range("J11").Value = "Gr. 70"
$ws.Range("L11").Value = "200 DEG. C"
$ws.Range("M11").Value = "14.00 kg/cm2 (g)"
$ws.Range("N11").Value = "150 DEG. C"
$ws.Range("O11").Value = "10.00 kg/cm2 (g)"

# Row 12
$ws.Range("E12").Value = "Tube"
$ws.Range("G12").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H12").Value = "Not Found"
$ws.Range("I12").Value = "SA-179"
$ws.Range("J12").ClearContents()
$ws.Range("L12").Value = "200 DEG. C"
$ws.Range("M12").Value = "14.00 kg/cm2 (g)"
$ws.Range("N12").Value = "150 DEG. C"
$ws.Range("O12").Value = "10.00 kg/cm2 (g)"

# Row 13
$ws.Range("E13").Value = "Tie Rods"
$ws.Range("G13").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H13").Value = "Carbon Steel"
$ws.Range("I13").Value = "SA-106"
$ws.Range("J13").Value = "Gr. B"
$ws.Range("L13").Value = "200 DEG. C"
$ws.Range("M13").Value = "14.00 kg/cm2 (g)"
$ws.Range("N13").Value = "150 DEG. C"
$ws.Range("O13").Value = "10.00 kg/cm2 (g)"

# Row 14
$ws.Range("E14").Value = "Baffle/Support Plate"
$ws.Range("G14").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H14").Value = "Carbon Steel"
$ws.Range("I14").Value = "SA-516"
$ws.Range("J14").Value = "Gr. 70"
$ws.Range("L14").Value = "200 DEG. C"
$ws.Range("M14").Value = "14.00 kg/cm2 (g)"
$ws.Range("N14").Value = "150 DEG. C"
$ws.Range("O14").Value = "10.00 kg/cm2 (g)"

# Row 15
$ws.Range("E15").Value = "Nozzle"
$ws.Range("G15").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H15").Value = "Carbon Steel"
$ws.Range("I15").Value = "SA-106"
$ws.Range("J15").Value = "Gr. B"
$ws.Range("L15").Value = "200 DEG. C"
$ws.Range("M15").Value = "14.00 kg/cm2 (g)"
$ws.Range("N15").Value = "150 DEG. C"
$ws.Range("O15").Value = "10.00 kg/cm2 (g)"

# Row 16
$ws.Range("E16").Value = "Saddle"
$ws.Range("G16").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H16").Value = "Carbon Steel"
$ws.Range("I16").Value = "SA-516"
$ws.Range("J16").Value = "Gr. 70"
$ws.Range("L16").Value = "200 DEG. C"
$ws.Range("M16").Value = "14.00 kg/cm2 (g)"
$ws.Range("N16").Value = "150 DEG. C"
$ws.Range("O16").Value = "10.00 kg/cm2 (g)"

# Row 17
$ws.Range("E17").Value = "Stiffening Rib"
$ws.Range("G17").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H17").Value = "Carbon Steel"
$ws.Range("I17").Value = "SA-516"
$ws.Range("J17").Value = "Gr. 70"
$ws.Range("L17").Value = "200 DEG. C"
$ws.Range("M17").Value = "14.00 kg/cm2 (g)"
$ws.Range("N17").Value = "150 DEG. C"
$ws.Range("O17").Value = "10.00 kg/cm2 (g)"

# Row 18
$ws.Range("E18").Value = "Stiffening Ring"
$ws.Range("G18").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H18").Value = "Carbon Steel"
$ws.Range("I18").Value = "SA-516"
$ws.Range("J18").Value = "Gr. 70"
$ws.Range("L18").Value = "200 DEG. C"
$ws.Range("M18").Value = "14.00 kg/cm2 (g)"
$ws.Range("N18").Value = "150 DEG. C"
$ws.Range("O18").Value = "10.00 kg/cm2 (g)"

# Row 19
$ws.Range("E19").Value = "Lifting Lug"
$ws.Range("G19").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H19").Value = "Carbon Steel"
$ws.Range("I19").Value = "SA-516"
$ws.Range("J19").Value = "Gr. 70"
$ws.Range("L19").Value = "200 DEG. C"
$ws.Range("M19").Value = "14.00 kg/cm2 (g)"
$ws.Range("N19").Value = "150 DEG. C"
$ws.Range("O19").Value = "10.00 kg/cm2 (g)"

# Row 20
$ws.Range("E20").Value = "Base Plate"
$ws.Range("G20").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H20").Value = "Carbon Steel"
$ws.Range("I20").Value = "SA-516"
$ws.Range("J20").Value = "Gr. 70"
$ws.Range("L20").Value = "200 DEG. C"
$ws.Range("M20").Value = "14.00 kg/cm2 (g)"
$ws.Range("N20").Value = "150 DEG. C"
$ws.Range("O20").Value = "10.00 kg/cm2 (g)"

# Row 21
$ws.Range("E21").Value = "Anchor Bolt"
$ws.Range("G21").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H21").Value = "Stainless Steel Bolting"
$ws.Range("I21").Value = "SA-193"
$ws.Range("J21").Value = "Gr. B7"
$ws.Range("L21").Value = "200 DEG. C"
$ws.Range("M21").Value = "14.00 kg/cm2 (g)"
$ws.Range("N21").Value = "150 DEG. C"
$ws.Range("O21").Value = "10.00 kg/cm2 (g)"

# --- New rows 22-29 ---
# Row 22
$ws.Range("E22").Value = "Hex. Nut"
$ws.Range("G22").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H22").Value = "Heavy Hex Nuts"
$ws.Range("I22").Value = "SA-194"
$ws.Range("J22").Value = "Gr. 2H"
$ws.Range("K22").Value = "N/A"
$ws.Range("L22").Value = "200 DEG. C"
$ws.Range("M22").Value = "14.00 kg/cm2 (g)"
$ws.Range("N22").Value = "150 DEG. C"
$ws.Range("O22").Value = "10.00 kg/cm2 (g)"

# Row 23
$ws.Range("E23").Value = "Washer"
$ws.Range("G23").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H23").Value = "Not Found"
$ws.Range("I23").Value = "SA-36"
$ws.Range("K23").Value = "N/A"
$ws.Range("L23").Value = "200 DEG. C"
$ws.Range("M23").Value = "14.00 kg/cm2 (g)"
$ws.Range("N23").Value = "150 DEG. C"
$ws.Range("O23").Value = "10.00 kg/cm2 (g)"

# Row 24
$ws.Range("E24").Value = "Name Plate"
$ws.Range("G24").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H24").Value = "Not Found"
$ws.Range("I24").Value = "SS-316"
$ws.Range("K24").Value = "N/A"
$ws.Range("L24").Value = "200 DEG. C"
$ws.Range("M24").Value = "14.00 kg/cm2 (g)"
$ws.Range("N24").Value = "150 DEG. C"
$ws.Range("O24").Value = "10.00 kg/cm2 (g)"

# Row 25
$ws.Range("E25").Value = "Stud Bolt for Flange"
$ws.Range("G25").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H25").Value = "Stainless Steel Bolting"
$ws.Range("I25").Value = "SA-193"
$ws.Range("J25").Value = "Gr. B7"
$ws.Range("K25").Value = "N/A"
$ws.Range("L25").Value = "200 DEG. C"
$ws.Range("M25").Value = "14.00 kg/cm2 (g)"
$ws.Range("N25").Value = "150 DEG. C"
$ws.Range("O25").Value = "10.00 kg/cm2 (g)"

# Row 26
$ws.Range("E26").Value = "Hex. Nut for Flange"
$ws.Range("G26").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H26").Value = "Heavy Hex Nuts"
$ws.Range("I26").Value = "SA-194"
$ws.Range("J26").Value = "Gr. 2H"
$ws.Range("K26").Value = "N/A"
$ws.Range("L26").Value = "200 DEG. C"
$ws.Range("M26").Value = "14.00 kg/cm2 (g)"
$ws.Range("N26").Value = "150 DEG. C"
$ws.Range("O26").Value = "10.00 kg/cm2 (g)"

# Row 27
$ws.Range("E27").Value = "Gasket"
$ws.Range("G27").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H27").Value = "Not Found"
$ws.Range("I27").Value = "Non Asbestos"
$ws.Range("K27").Value = "N/A"
$ws.Range("L27").Value = "200 DEG. C"
$ws.Range("M27").Value = "14.00 kg/cm2 (g)"
$ws.Range("N27").Value = "150 DEG. C"
$ws.Range("O27").Value = "10.00 kg/cm2 (g)"

# Row 28
$ws.Range("E28").Value = "Vent / Drain"
$ws.Range("G28").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H28").Value = "Carbon Steel"
$ws.Range("I28").Value = "SA-106"
$ws.Range("J28").Value = "Gr. B"
$ws.Range("K28").Value = "N/A"
$ws.Range("L28").Value = "200 DEG. C"
$ws.Range("M28").Value = "14.00 kg/cm2 (g)"
$ws.Range("N28").Value = "150 DEG. C"
$ws.Range("O28").Value = "10.00 kg/cm2 (g)"

# Row 29
$ws.Range("E29").Value = "Tube Support"
$ws.Range("G29").Value = "METHANOL (SHELL) / WATER (TUBE)"
$ws.Range("H29").Value = "Carbon Steel"
$ws.Range("I29").Value = "SA-516"
$ws.Range("J29").Value = "Gr. 70"
$ws.Range("K29").Value = "N/A"
$ws.Range("L29").Value = "200 DEG. C"
$ws.Range("M29").Value = "14.00 kg/cm2 (g)"
$ws.Range("N29").Value = "150 DEG. C"
$ws.Range("O29").Value = "10.00 kg/cm2 (g)"

# --- Merge cell range updates (A,B,C,D columns: 8:21 -> 8:29) ---
$ws.Range("A8:A21").UnMerge()
$ws.Range("A8:A29").Merge()
$ws.Range("B8:B21").UnMerge()
$ws.Range("B8:B29").Merge()
$ws.Range("C8:C21").UnMerge()
$ws.Range("C8:C29").Merge()
$ws.Range("D8:D21").UnMerge()
$ws.Range("D8:D29").Merge()

Write-Output "done"